$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "901/FES "
$ws.Range("C2").Value = "J207703"
$ws.Range("D2").Value = "ACHENGLI LAILA"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 5000
$ws.Range("J2").Value = 500
$ws.Range("O2").Value = 4500

# Row 3 updates
$ws.Range("A3").Value = "901/LF/FES "
$ws.Range("B3").Value = "Logement de fonction"
$ws.Range("C3").Value = "BJ36877"
$ws.Range("D3").Value = "CHARIJI ABDELLAH"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 6000
$ws.Range("J3").Value = 600
$ws.Range("O3").Value = 5400

# Row 4 updates - becomes blank-ish row with spaces
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = 11000
$ws.Range("J4").Value = 1100
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 9900

# Delete row 5 entirely
$ws.Rows("5").Delete()
